# Edit and re export DD Gantt
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Task name text update (B4) ---
$ws.Range("B4").Value = "High level system view definition"

# --- Reorder tasks: swap row 10 <-> row 11 (Task Name / Resource) ---
$b10 = $ws.Range("B10").Text
$c10 = $ws.Range("C10").Text
$b11 = $ws.Range("B11").Text
$c11 = $ws.Range("C11").Text
$ws.Range("B10").Value = $b11
$ws.Range("C10").Value = $c11
$ws.Range("B11").Value = $b10
$ws.Range("C11").Value = $c10

# --- Reorder tasks: swap row 13 <-> row 14 (Task Name / Resource) ---
$b13 = $ws.Range("B13").Text
$c13 = $ws.Range("C13").Text
$b14 = $ws.Range("B14").Text
$c14 = $ws.Range("C14").Text
$ws.Range("B13").Value = $b14
$ws.Range("C13").Value = $c14
$ws.Range("B14").Value = $b13
$ws.Range("C14").Value = $c13

# --- New task note on row 13 ---
$ws.Range("D13").Value = "Moreno will start working on this activity in date 29/11/2016"

# --- Updated schedule (Start Date / Finish Date) ---
$ws.Range("G7").Value = 42691.666666666664

$ws.Range("F8").Value = 42692.333333333336
$ws.Range("G8").Value = 42695.666666666664

$ws.Range("F9").Value = 42696.333333333336
$ws.Range("G9").Value = 42697.666666666664

$ws.Range("F10").Value = 42698.333333333336
$ws.Range("G10").Value = 42698.666666666664

$ws.Range("F11").Value = 42698.333333333336
$ws.Range("G11").Value = 42699.666666666664

$ws.Range("F12").Value = 42698.333333333336
$ws.Range("G12").Value = 42699.666666666664

$ws.Range("F13").Value = 42699.333333333336
$ws.Range("G13").Value = 42704.666666666664

$ws.Range("F14").Value = 42702.333333333336
$ws.Range("G14").Value = 42702.666666666664

$ws.Range("F15").Value = 42703.333333333336
$ws.Range("G15").Value = 42704.666666666664

$ws.Range("F16").Value = 42705.333333333336
$ws.Range("G16").Value = 42705.666666666664

$ws.Range("F17").Value = 42705.333333333336
$ws.Range("G17").Value = 42705.666666666664

$ws.Range("F18").Value = 42706.333333333336
$ws.Range("G18").Value = 42706.666666666664

# --- Move the visible selection to match the saved view state ---
$ws.Range("B27").Select()
